# Insert a new data row at row 381 (pushes the existing rows 381..479 down
# to 382..480, matching the diff's dimension change A1:R479 -> A1:R480) and
# populate it with the new "Apio" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(381).Insert()

$ws.Range("A381").Value = 3
$ws.Range("B381").Value = 'Femacal de La Calera'
$ws.Range("C381").Value = 'Coquimbo'
$ws.Range("D381").Value = 44855
$ws.Range("E381").Value = 5
$ws.Range("F381").Value = 100112017
$ws.Range("G381").Value = 'Apio'
$ws.Range("H381").Value = 'Americana (o)'
$ws.Range("I381").Value = 'Primera'
$ws.Range("J381").Value = 230
$ws.Range("K381").Value = 9000
$ws.Range("L381").Value = 9500
$ws.Range("M381").Value = 9239
$ws.Range("N381").Value = '$/docena de matas'
$ws.Range("O381").Value = 'Pan de Azúcar'
$ws.Range("P381").Value = 1540
$ws.Range("Q381").Value = 6
$ws.Range("R381").Value = 'Hortaliza'
